$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # The new C4/E4/G4 "% Dev" cells need the same formatting (style index 2:
    # centered alignment) as the rest of row 4 - copy it from the existing B4
    # cell rather than leaving the new cells with the default style.
    $ws.Range("B4").Copy() | Out-Null
    $ws.Range("C4").PasteSpecial(-4122) | Out-Null
    $ws.Range("E4").PasteSpecial(-4122) | Out-Null
    $ws.Range("G4").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    # Fill in the computed sublimation-energy values that the input parser
    # had previously zeroed out.
    $ws.Range("B4").Value = 109.7451299964792
    $ws.Range("C4").Value = 0.5840038070334835
    $ws.Range("D4").Value = 109.7451299964792
    $ws.Range("E4").Value = 0.5840038070334835
    $ws.Range("F4").Value = 109.7451299964792
    $ws.Range("G4").Value = 0.5840038070334835
}
